$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 874 entirely (the "行くな / ماتروحش" post), which shifts all
# subsequent rows up by one and shrinks the used range from A1:C888 to A1:C887.
$ws.Rows.Item(874).Delete()
